# fix VideoView cannot play in 2.3.6 / auto mode S39H cannot play videoview
#
# Adds two new rows (No. 24 and No. 25) to the "bug" worksheet's bug list,
# matching the formatting of the row immediately above them (row 24 / No. 23),
# and gives the last (new) cell a highlighted, word-wrapped style for its
# long multi-line description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone formatting from the last existing data row (row 24) into the
#        row being completed (row 25) and the brand-new row (row 26). -------
$ws.Range("B24:H24").Copy()
$ws.Range("B25:H25").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B26:H26").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Row 25 (bug No. 24) --------------------------------------------------
$ws.Range("B25").Value = 24
$ws.Range("C25").Value = "N/A"
$ws.Range("D25").Value = "android"
$ws.Range("E25").Value = 20150123
$ws.Range("F25").Value = "2.3.6手机 变形金刚 截图很容易失败"
$ws.Range("G25").Value = "tracking"

# --- 3. Row 26 (bug No. 25, new row) ----------------------------------------
$ws.Range("B26").Value = 25
$ws.Range("C26").Value = "N/A"
$ws.Range("D26").Value = "IOS"
$ws.Range("E26").Value = 20150116

$newLine = [char]10
$f26Text = "[铁臂阿童木][Astro.Boy_48k.mp3 声音不对" + $newLine + "国际潜水小姐大赛.冠军专访.m2t" + $newLine + "D:\Archive\media\audio\邓紫棋 - 泡沫.flac" + $newLine + "D:\Archive\media\audio\陈慧娴-飘雪.ape"

# F26 gets its own highlighted / wrapped style (distinct from the rest of the row)
$ws.Range("F26").Interior.Color = 65535
$ws.Range("F26").HorizontalAlignment = -4108
$ws.Range("F26").VerticalAlignment = -4108
$ws.Range("F26").WrapText = $true
$ws.Range("F26").Value = $f26Text

$ws.Range("G26").Value = "tracking"

# Row 26 is tall enough to show the 4-line description.
$ws.Rows.Item(26).RowHeight = 60

# --- 4. Selection mirrors where the edit ended up. --------------------------
$ws.Range("F26").Select()
